$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header at H1, matching the style of the other headers (copy format from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values for rows 2-10 (era data updated alongside)
$saveValues = @(1, 0, 0, 0, 0, 0, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
